# fall 24 week 5 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 10.61
$ws.Range("H4").Value = 9.140000000000001

$ws.Range("D5").Value = 9.390000000000001
$ws.Range("F5").Value = 10.15
$ws.Range("G5").Value = 9.68
$ws.Range("J5").Value = 7.56

$ws.Range("E6").Value = 9.85
$ws.Range("G6").Value = 10.47
$ws.Range("H6").Value = 10.43

$ws.Range("E7").Value = 10.32
$ws.Range("F7").Value = 9.529999999999999

$ws.Range("D8").Value = 10.86
$ws.Range("F8").Value = 9.57

$ws.Range("E10").Value = 12.44
